$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.750.78'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '2.619.06'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '515.75'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.43'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = '2.633.99'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  +4.14%  '
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('D14').Value = '3.075.98'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '60.756.74'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('D18').Value = '2.624.77'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '357.74'
$ws.Range('E20').Value = '  +4.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.66'
$ws.Range('E21').Value = '  +2.69%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.89'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = '2.735.49'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.49'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '151.54'
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.93'
$ws.Range('E35').Value = '  +3.95%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  +6.86%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.49'
$ws.Range('E39').Value = '  +1.63%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.855'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.40'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.76'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '292.91'
$ws.Range('E43').Value = '  -4.55%  '
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.623'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0558'
$ws.Range('E46').Value = '  -2.18%  '
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.81'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.00'
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.30'
$ws.Range('E51').Value = '  +0.35%  '
